# Weekly fruit/vegetable price update: swap the data recorded for the two
# "Primera" quality Níspero entries (rows 2 & 4) and the two remaining
# entries (rows 5 & 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -> becomes what row 4 used to hold ---
$ws.Range("D2").Value = 44488
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("Q2").Value = '$/bandeja 5 kilos'
$ws.Range("R2").Value = 'La Ligua'
$ws.Range("S2").Value = 2400
$ws.Range("T2").Value = 5

# --- Row 4 -> becomes what row 2 used to hold ---
$ws.Range("D4").Value = 44496
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("Q4").Value = '$/bandeja 10 kilos'
$ws.Range("R4").Value = 'Provincia de Quillota'
$ws.Range("S4").Value = 2800
$ws.Range("T4").Value = 10

# --- Row 5 -> becomes what row 6 used to hold ---
$ws.Range("D5").Value = 44466
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 11000
$ws.Range("Q5").Value = '$/bandeja 5 kilos'
$ws.Range("S5").Value = 2200
$ws.Range("T5").Value = 5

# --- Row 6 -> becomes what row 5 used to hold ---
$ws.Range("D6").Value = 44166
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("S6").Value = 667
$ws.Range("T6").Value = 18
